$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting (bold, border, centered alignment) from the existing
# header cell H1 onto the two new header cells so they reuse the same style.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Re-assert the values (PasteSpecial only touches formats, but keep this to
# be safe in case paste affected contents).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values for columns I (I0) and J (IF), rows 2-10
$iValues = @(7, 9, 5, 8, 8, 7, 5, 4, 7)
$jValues = @(9, 9, 7, 8, 9, 7, 6, 6, 8)

for ($r = 0; $r -lt 9; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$r]
    $ws.Cells.Item($row, 10).Value = $jValues[$r]
}
